# Realestate Update resale numbers 2024-01-02 12:27
# Append a new data row (row 6) to the single worksheet, mirroring the
# existing rows: columns A-D are text (date/time/weekday/week-code,
# the "00" week code must keep its leading zero), columns E-T are plain
# numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the text columns to be stored as text rather than being
# auto-parsed into a date serial / number (Excel would otherwise turn
# "2024-01-02" into a date and "00" into 0).
$ws.Range("A6:D6").NumberFormat = "@"

$ws.Range("A6").Value = "2024-01-02"
$ws.Range("B6").Value = "12:27:14"
$ws.Range("C6").Value = "Tuesday"
$ws.Range("D6").Value = "00"

# Drop the temporary "Text" number format again so the new cells end up
# styled the same (no explicit style) as the rest of the sheet - only
# the stored value/type is affected by ClearFormats, not the text we
# already committed above.
$ws.Range("A6:D6").ClearFormats()

$ws.Range("E6").Value = 140184
$ws.Range("F6").Value = 142929
$ws.Range("G6").Value = 171456
$ws.Range("H6").Value = 145698
$ws.Range("I6").Value = -1
$ws.Range("J6").Value = 116840
$ws.Range("K6").Value = 223711
$ws.Range("L6").Value = 247756
$ws.Range("M6").Value = 183519
$ws.Range("N6").Value = 109856
$ws.Range("O6").Value = 39630
$ws.Range("P6").Value = 30657
$ws.Range("Q6").Value = 71726
$ws.Range("R6").Value = -1
$ws.Range("S6").Value = 41410
$ws.Range("T6").Value = -1
